$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# Edit 1: the "מגישות:" heading paragraph currently carries an explicit
# paragraph-mark formatting of <w:rPr><w:rtl/></w:rPr> inside its <w:pPr>.
# The edit drops that paragraph-mark formatting entirely (the run itself,
# with its own bold/underline/rtl formatting, is left untouched).
# ---------------------------------------------------------------------------
$targetText = "מגישות:"
$heading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]10, [char]7)
    if ($t -eq $targetText) {
        $heading = $p
        break
    }
}

if ($heading -ne $null) {
    $headingXml = '<w:p xmlns:w="' + $wNs + '">' +
                    '<w:r>' +
                      '<w:rPr><w:rFonts w:hint="cs"/><w:b/><w:bCs/><w:u w:val="single"/><w:rtl/></w:rPr>' +
                      '<w:t>' + $targetText + '</w:t>' +
                    '</w:r>' +
                  '</w:p>'
    $heading.Range.InsertXML($headingXml)
}

# ---------------------------------------------------------------------------
# Edit 2: the final (empty) paragraph of the document - which had
# <w:ind w:firstLine="720"/> plus <w:rPr><w:rFonts w:hint="cs"/><w:rtl/></w:rPr>
# in its <w:pPr> - is turned into two paragraphs:
#   1) an empty paragraph whose paragraph mark is bold/underlined/rtl
#      (no first-line indent any more)
#   2) a new paragraph containing the run "הערה: " (bold, underlined, rtl)
# ---------------------------------------------------------------------------
$last = $d.Paragraphs.Item($d.Paragraphs.Count)

$newXml = '<w:p xmlns:w="' + $wNs + '">' +
            '<w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:rtl/></w:rPr></w:pPr>' +
          '</w:p>' +
          '<w:p xmlns:w="' + $wNs + '">' +
            '<w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr>' +
            '<w:r>' +
              '<w:rPr><w:rFonts w:hint="cs"/><w:b/><w:bCs/><w:u w:val="single"/><w:rtl/></w:rPr>' +
              '<w:t xml:space="preserve">הערה: </w:t>' +
            '</w:r>' +
          '</w:p>'

$last.Range.InsertXML($newXml)
